$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header D1: "Low EGU limit" -> "Shelve policy"
$ws.Range("D1").Value = "Shelve policy"

# A2 and A3: "AI" -> "DI"
$ws.Range("A2").Value = "DI"
$ws.Range("A3").Value = "DI"
